# Apply crypto price/volume update (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.120.12'
$ws.Range('E2').Value = '  -1.60%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.996.06'
$ws.Range('E3').Value = '  -0.58%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.015'
$ws.Range('E4').Value = '  +0.62%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '330.79'
$ws.Range('E5').Value = '  -0.22%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.013'
$ws.Range('E6').Value = '  +0.59%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4976'
$ws.Range('E7').Value = '  -1.02%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4190'
$ws.Range('E8').Value = '  -1.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '54.85'
$ws.Range('E9').Value = '  +1.86%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.08836'
$ws.Range('E10').Value = '  -4.00%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.094'
$ws.Range('E11').Value = '  -2.99%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.93'
$ws.Range('E12').Value = '  -2.47%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.993.43'
$ws.Range('E13').Value = '  +1.93%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '7.988'
$ws.Range('E14').Value = '  -1.55%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.427'
$ws.Range('E15').Value = '  -1.83%  '
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '92.52'
$ws.Range('E17').Value = '  -3.46%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001104'
$ws.Range('E18').Value = '  -1.99%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06738'
$ws.Range('E19').Value = '  +1.11%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '19.48'
$ws.Range('E20').Value = '  -1.84%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.013'
$ws.Range('E21').Value = '  +0.72%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.977'
$ws.Range('E22').Value = '  -0.31%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.156.30'
$ws.Range('E23').Value = '  -1.53%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.98'
$ws.Range('E24').Value = '  -0.10%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.293'
$ws.Range('E25').Value = '  +0.54%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.257.47'
$ws.Range('E26').Value = '  +3.97%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '20.84'
$ws.Range('E27').Value = '  +0.18%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '157.00'
$ws.Range('E28').Value = '  -1.31%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.292'
$ws.Range('E29').Value = '  -4.20%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.254'
$ws.Range('E30').Value = '  -3.89%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '127.15'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.044'
$ws.Range('E32').Value = '  -1.17%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09864'
$ws.Range('E33').Value = '  -1.10%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.531'
$ws.Range('E34').Value = '  -4.10%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.821'
$ws.Range('E35').Value = '  -0.95%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.737'
$ws.Range('E36').Value = '  -1.34%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02420'
$ws.Range('E37').Value = '  -2.35%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '9.146'
$ws.Range('E38').Value = '  -5.26%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.311'
$ws.Range('E39').Value = '  -0.07%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06384'
$ws.Range('E40').Value = '  -0.02%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6485'
$ws.Range('E41').Value = '  -1.37%  '
$ws.Range('E42').Value = '  -1.79%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1981'
$ws.Range('E43').Value = '  -4.74%  '
$ws.Range('E44').Value = '  +0.55%  '
$ws.Range('B45').Value = 'Decentraland'
$ws.Range('C45').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.6201'
$ws.Range('E45').Value = '  -2.48%  '
$ws.Range('B46').Value = 'WEMIXTOKEN'
$ws.Range('C46').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.356'
$ws.Range('E46').Value = '  +4.13%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '13.41'
$ws.Range('E47').Value = '  -0.10%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.172'
$ws.Range('E48').Value = '  -2.02%  '
$ws.Range('B49').Value = 'BabyDogeCoin'
$ws.Range('C49').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.00000000344'
$ws.Range('E49').Value = '  +4.86%  '
$ws.Range('B50').Value = 'PancakeSwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.494'
$ws.Range('E50').Value = '  -1.08%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.180'
$ws.Range('E51').Value = '  +7.98%  '
